$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 677.1111
$ws.Range("I4").Value = 510.5
$ws.Range("K4").Value = 510.5
$ws.Range("M4").Value = -396.5
$ws.Range("H17").Value = 14551.889
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 14551.889
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 43655.667
$ws.Range("N17").Value = -43991.667
$ws.Range("H31").Value = 1999
$ws.Range("I31").Value = 1999
$ws.Range("K31").Value = 5997
$ws.Range("M31").Value = -5767
$ws.Range("H41").Value = 3344.9524
$ws.Range("I41").Value = 3519.4707
$ws.Range("J41").Value = 2603.25
$ws.Range("K41").Value = 3519.4707
$ws.Range("L41").Value = 2603.25
$ws.Range("M41").Value = -3079.4707
$ws.Range("N41").Value = -3483.25
$ws.Range("H43").Value = 2266.6667
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("H53").Value = 1201.1923
$ws.Range("I53").Value = 709.5
$ws.Range("K53").Value = 709.5
$ws.Range("M53").Value = -72.5
$ws.Range("H58").Value = 2526
$ws.Range("J58").Value = 4663.3335
$ws.Range("L58").Value = 13990.0005
$ws.Range("N58").Value = -14290.0005
$ws.Range("H62").Value = 6196.273
$ws.Range("I62").Value = 5776
$ws.Range("K62").Value = 5776
$ws.Range("M62").Value = -5152
$ws.Range("H65").Value = 6196.273
$ws.Range("I65").Value = 5776
$ws.Range("K65").Value = 28880
$ws.Range("M65").Value = -25760
$ws.Range("H68").Value = 120999
$ws.Range("J68").Value = 120999
$ws.Range("L68").Value = 120999
$ws.Range("N68").Value = -122497
$ws.Range("H70").Value = 41668210
$ws.Range("I70").Value = 1959.3334
$ws.Range("J70").Value = 83334456
$ws.Range("K70").Value = 5878.0002
$ws.Range("L70").Value = 250003368
$ws.Range("M70").Value = -5608.0002
$ws.Range("N70").Value = -250003908
$ws.Range("H71").Value = 120999
$ws.Range("J71").Value = 120999
$ws.Range("L71").Value = 362997
$ws.Range("N71").Value = -370485
$ws.Range("H73").Value = 41668210
$ws.Range("I73").Value = 1959.3334
$ws.Range("J73").Value = 83334456
$ws.Range("K73").Value = 5878.0002
$ws.Range("L73").Value = 250003368
$ws.Range("M73").Value = -4942.0002
$ws.Range("N73").Value = -250005240
$ws.Range("H76").Value = 6449.5
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("H79").Value = 6449.5
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("H86").Value = 75004870
$ws.Range("I86").Value = 55560784
$ws.Range("J86").Value = 104171010
$ws.Range("K86").Value = 55560784
$ws.Range("L86").Value = 104171010
$ws.Range("M86").Value = -55559661
$ws.Range("N86").Value = -104173256
$ws.Range("H89").Value = 75004870
$ws.Range("I89").Value = 55560784
$ws.Range("J89").Value = 104171010
$ws.Range("K89").Value = 277803920
$ws.Range("L89").Value = 520855050
$ws.Range("M89").Value = -277798304
$ws.Range("N89").Value = -520866282
$ws.Range("H92").Value = 21741264
$ws.Range("I92").Value = 22728140
$ws.Range("K92").Value = 22728140
$ws.Range("M92").Value = -22726892
$ws.Range("H98").Value = 2230.138
$ws.Range("I98").Value = 2230.138
$ws.Range("K98").Value = 2230.138
$ws.Range("M98").Value = -732.1379999999999
$ws.Range("H100").Value = 3380.6667
$ws.Range("I100").Value = 1656.8
$ws.Range("K100").Value = 1656.8
$ws.Range("M100").Value = -1115.8
$ws.Range("H106").Value = 33335332
$ws.Range("I106").Value = 33335332
$ws.Range("K106").Value = 33335332
$ws.Range("M106").Value = -33334701
$ws.Range("H112").Value = 3684.06
$ws.Range("J112").Value = 3756.3125
$ws.Range("L112").Value = 11268.9375
$ws.Range("N112").Value = -13484.9375
$ws.Range("H122").Value = 2230.138
$ws.Range("I122").Value = 2230.138
$ws.Range("K122").Value = 6690.414
$ws.Range("M122").Value = -4240.414
$ws.Range("H129").Value = 2763.6572
$ws.Range("I129").Value = 912.0833
$ws.Range("K129").Value = 2736.2499
$ws.Range("M129").Value = 2263.7501
$ws.Range("H132").Value = 3167.9636
$ws.Range("I132").Value = 3035.4695
$ws.Range("K132").Value = 9106.408500000001
$ws.Range("M132").Value = -6576.408500000001
$ws.Range("H135").Value = 1169.1936
$ws.Range("J135").Value = 1999.6666
$ws.Range("L135").Value = 17996.9994
$ws.Range("N135").Value = -23066.9994
$ws.Range("H137").Value = 3207.0908
$ws.Range("I137").Value = 3831.5
$ws.Range("J137").Value = 2457.8
$ws.Range("K137").Value = 11494.5
$ws.Range("L137").Value = 7373.400000000001
$ws.Range("M137").Value = -8944.5
$ws.Range("N137").Value = -12473.4
$ws.Range("H138").Value = 2920.94
$ws.Range("I138").Value = 2145.2
$ws.Range("J138").Value = 3179.52
$ws.Range("K138").Value = 6435.599999999999
$ws.Range("L138").Value = 9538.56
$ws.Range("M138").Value = -1295.599999999999
$ws.Range("N138").Value = -19818.56
$ws.Range("H140").Value = 272137.84
$ws.Range("I140").Value = 239985
$ws.Range("J140").Value = 277496.66
$ws.Range("K140").Value = 239985
$ws.Range("L140").Value = 277496.66
$ws.Range("M140").Value = -234805
$ws.Range("N140").Value = -287856.66
$ws.Range("M43").ClearContents()
$ws.Range("N76").ClearContents()
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1756.5555
$ws.Range("I2").Value = 1711.9333
$ws.Range("K2").Value = 1711.9333
$ws.Range("M2").Value = -1598.9333
$ws.Range("H32").Value = 8208.83
$ws.Range("I32").Value = 5319.7363
$ws.Range("J32").Value = 15637.929
$ws.Range("K32").Value = 5319.7363
$ws.Range("L32").Value = 15637.929
$ws.Range("M32").Value = -5032.7363
$ws.Range("N32").Value = -16211.929
$ws.Range("H37").Value = 14052.667
$ws.Range("I37").Value = 14052.667
$ws.Range("K37").Value = 14052.667
$ws.Range("M37").Value = -13779.667
$ws.Range("H45").Value = 5964.4
$ws.Range("I45").Value = 5966
$ws.Range("J45").Value = 5950
$ws.Range("K45").Value = 5966
$ws.Range("L45").Value = 5950
$ws.Range("M45").Value = -5589
$ws.Range("N45").Value = -6704
$ws.Range("H61").Value = 3482.1428
$ws.Range("I61").Value = 3063.389
$ws.Range("K61").Value = 3063.389
$ws.Range("M61").Value = -2851.389
$ws.Range("H63").Value = 159533.89
$ws.Range("I63").Value = 3001.6667
$ws.Range("K63").Value = 3001.6667
$ws.Range("M63").Value = -2315.6667
$ws.Range("H66").Value = 159533.89
$ws.Range("I66").Value = 3001.6667
$ws.Range("K66").Value = 15008.3335
$ws.Range("M66").Value = -11576.3335
$ws.Range("H74").Value = 2963.762
$ws.Range("I74").Value = 2757.9285
$ws.Range("J74").Value = 3375.4285
$ws.Range("K74").Value = 2757.9285
$ws.Range("L74").Value = 3375.4285
$ws.Range("M74").Value = -1883.9285
$ws.Range("N74").Value = -5123.4285
$ws.Range("H77").Value = 2963.762
$ws.Range("I77").Value = 2757.9285
$ws.Range("J77").Value = 3375.4285
$ws.Range("K77").Value = 13789.6425
$ws.Range("L77").Value = 16877.1425
$ws.Range("M77").Value = -9421.6425
$ws.Range("N77").Value = -25613.1425
$ws.Range("H88").Value = 1646
$ws.Range("J88").Value = 2005.5
$ws.Range("L88").Value = 2005.5
$ws.Range("N88").Value = -2817.5
$ws.Range("H91").Value = 1646
$ws.Range("J91").Value = 2005.5
$ws.Range("L91").Value = 2005.5
$ws.Range("N91").Value = -4813.5
$ws.Range("H97").Value = 946.8571
$ws.Range("I97").Value = 611.0909
$ws.Range("J97").Value = 2178
$ws.Range("K97").Value = 611.0909
$ws.Range("L97").Value = 2178
$ws.Range("M97").Value = -115.0909
$ws.Range("N97").Value = -3170
$ws.Range("H110").Value = 1577.1818
$ws.Range("I110").Value = 1195.1578
$ws.Range("K110").Value = 1195.1578
$ws.Range("M110").Value = 849.8422
$ws.Range("H116").Value = 1756.5555
$ws.Range("I116").Value = 1711.9333
$ws.Range("K116").Value = 1711.9333
$ws.Range("M116").Value = 582.0667000000001
$ws.Range("H122").Value = 5962.0303
$ws.Range("I122").Value = 4654.2173
$ws.Range("K122").Value = 13962.6519
$ws.Range("M122").Value = -11512.6519
$ws.Range("H132").Value = 117847.66
$ws.Range("I132").Value = 163738.83
$ws.Range("K132").Value = 491216.49
$ws.Range("M132").Value = -488686.49
$ws.Range("H136").Value = 3482.1428
$ws.Range("I136").Value = 3063.389
$ws.Range("K136").Value = 9190.167000000001
$ws.Range("M136").Value = -6640.167000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1756.5555
$ws.Range("I3").Value = 1711.9333
$ws.Range("K3").Value = 1711.9333
$ws.Range("M3").Value = -1597.9333
$ws.Range("H36").Value = 866
$ws.Range("I36").Value = 866
$ws.Range("K36").Value = 866
$ws.Range("M36").Value = -332
$ws.Range("H86").Value = 3996.5
$ws.Range("I86").Value = 3996.5
$ws.Range("K86").Value = 3996.5
$ws.Range("M86").Value = -2873.5
$ws.Range("H89").Value = 3996.5
$ws.Range("I89").Value = 3996.5
$ws.Range("K89").Value = 19982.5
$ws.Range("M89").Value = -14366.5
$ws.Range("H93").Value = 62025.5
$ws.Range("J93").Value = 62025.5
$ws.Range("L93").Value = 62025.5
$ws.Range("N93").Value = -65769.5
$ws.Range("H99").Value = 2073.35
$ws.Range("I99").Value = 1694.1428
$ws.Range("K99").Value = 1694.1428
$ws.Range("M99").Value = -196.1428000000001
$ws.Range("H132").Value = 126490
$ws.Range("J132").Value = 126490
$ws.Range("L132").Value = 126490
$ws.Range("N132").Value = -136610
$ws.Range("H133").Value = 20709
$ws.Range("I133").Value = 20709
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 20709
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -15649
$ws.Range("H134").Value = 1663153.6
$ws.Range("I134").Value = 2233600.5
$ws.Range("K134").Value = 6700801.5
$ws.Range("M134").Value = -6698266.5
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1991.0588
$ws.Range("I16").Value = 1703.4286
$ws.Range("J16").Value = 3333.3333
$ws.Range("K16").Value = 1703.4286
$ws.Range("L16").Value = 3333.3333
$ws.Range("M16").Value = -1416.4286
$ws.Range("N16").Value = -3907.3333
$ws.Range("H31").Value = 9801.483
$ws.Range("I31").Value = 6296.636
$ws.Range("K31").Value = 6296.636
$ws.Range("M31").Value = -6001.636
$ws.Range("H34").Value = 9801.483
$ws.Range("I34").Value = 6296.636
$ws.Range("K34").Value = 6296.636
$ws.Range("M34").Value = -6094.636
$ws.Range("H45").Value = 20500
$ws.Range("I45").Value = 35000
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 35000
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = -34407
$ws.Range("N45").Value = -7186
$ws.Range("H58").Value = 3341.2222
$ws.Range("I58").Value = 2997.6428
$ws.Range("J58").Value = 4543.75
$ws.Range("K58").Value = 2997.6428
$ws.Range("L58").Value = 4543.75
$ws.Range("M58").Value = -2794.6428
$ws.Range("N58").Value = -4949.75
$ws.Range("H62").Value = 4413.875
$ws.Range("I62").Value = 3971.5
$ws.Range("J62").Value = 4856.25
$ws.Range("K62").Value = 3971.5
$ws.Range("L62").Value = 4856.25
$ws.Range("M62").Value = -3347.5
$ws.Range("N62").Value = -6104.25
$ws.Range("H65").Value = 4413.875
$ws.Range("I65").Value = 3971.5
$ws.Range("J65").Value = 4856.25
$ws.Range("K65").Value = 19857.5
$ws.Range("L65").Value = 24281.25
$ws.Range("M65").Value = -16737.5
$ws.Range("N65").Value = -30521.25
$ws.Range("H75").Value = 113998.5
$ws.Range("J75").Value = 113998.5
$ws.Range("L75").Value = 113998.5
$ws.Range("N75").Value = -115994.5
$ws.Range("H78").Value = 113998.5
$ws.Range("J78").Value = 113998.5
$ws.Range("L78").Value = 341995.5
$ws.Range("N78").Value = -351979.5
$ws.Range("H86").Value = 3842
$ws.Range("I86").Value = 3870.4
$ws.Range("K86").Value = 3870.4
$ws.Range("M86").Value = -2747.4
$ws.Range("H88").Value = 34749.5
$ws.Range("J88").Value = 37999.668
$ws.Range("L88").Value = 37999.668
$ws.Range("N88").Value = -38811.668
$ws.Range("H89").Value = 3842
$ws.Range("I89").Value = 3870.4
$ws.Range("K89").Value = 19352
$ws.Range("M89").Value = -13736
$ws.Range("H91").Value = 34749.5
$ws.Range("J91").Value = 37999.668
$ws.Range("L91").Value = 37999.668
$ws.Range("N91").Value = -40807.668
$ws.Range("H94").Value = 1119.1818
$ws.Range("I94").Value = 1064.6
$ws.Range("K94").Value = 1064.6
$ws.Range("M94").Value = -613.5999999999999
$ws.Range("H95").Value = 40999.5
$ws.Range("J95").Value = 40999.5
$ws.Range("L95").Value = 40999.5
$ws.Range("N95").Value = -46491.5
$ws.Range("H103").Value = 57500.25
$ws.Range("I103").Value = 14999.75
$ws.Range("J103").Value = 100000.75
$ws.Range("K103").Value = 14999.75
$ws.Range("L103").Value = 100000.75
$ws.Range("M103").Value = -13827.75
$ws.Range("N103").Value = -102344.75
$ws.Range("H110").Value = 48995
$ws.Range("J110").Value = 37326.668
$ws.Range("L110").Value = 37326.668
$ws.Range("N110").Value = -45506.668
$ws.Range("H113").Value = 1991.0588
$ws.Range("I113").Value = 1703.4286
$ws.Range("J113").Value = 3333.3333
$ws.Range("K113").Value = 1703.4286
$ws.Range("L113").Value = 3333.3333
$ws.Range("M113").Value = 466.5714
$ws.Range("N113").Value = -7673.3333
$ws.Range("H132").Value = 1428.1428
$ws.Range("I132").Value = 1428.1428
$ws.Range("K132").Value = 4284.428400000001
$ws.Range("M132").Value = -1754.428400000001
$ws.Range("H134").Value = 1971.1945
$ws.Range("I134").Value = 1370.3704
$ws.Range("K134").Value = 4111.1112
$ws.Range("M134").Value = -1576.1112
$ws.Range("H136").Value = 3341.2222
$ws.Range("I136").Value = 2997.6428
$ws.Range("J136").Value = 4543.75
$ws.Range("K136").Value = 8992.9284
$ws.Range("L136").Value = 13631.25
$ws.Range("M136").Value = -6442.928400000001
$ws.Range("N136").Value = -18731.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 11304.4
$ws.Range("I14").Value = 11304.4
$ws.Range("K14").Value = 33913.2
$ws.Range("M14").Value = -33740.2
$ws.Range("H17").Value = 1296.6471
$ws.Range("J17").Value = 728
$ws.Range("L17").Value = 2184
$ws.Range("N17").Value = -2522
$ws.Range("H39").Value = 4972.8076
$ws.Range("J39").Value = 4972.8076
$ws.Range("L39").Value = 14918.4228
$ws.Range("N39").Value = -15506.4228
$ws.Range("H55").Value = 2100.158
$ws.Range("J55").Value = 2713.1538
$ws.Range("L55").Value = 8139.4614
$ws.Range("N55").Value = -8493.4614
$ws.Range("H62").Value = 3496.3333
$ws.Range("J62").Value = 3994.5
$ws.Range("L62").Value = 11983.5
$ws.Range("N62").Value = -13355.5
$ws.Range("H65").Value = 3496.3333
$ws.Range("J65").Value = 3994.5
$ws.Range("L65").Value = 35950.5
$ws.Range("N65").Value = -42814.5
$ws.Range("H104").Value = 7600.1113
$ws.Range("J104").Value = 8729.167
$ws.Range("L104").Value = 26187.501
$ws.Range("N104").Value = -31429.501
$ws.Range("H107").Value = 876.5
$ws.Range("J107").Value = 853.97296
$ws.Range("L107").Value = 2561.91888
$ws.Range("N107").Value = -6401.918879999999
$ws.Range("H129").Value = 1667.909
$ws.Range("I129").Value = 801.06665
$ws.Range("J129").Value = 3525.4285
$ws.Range("K129").Value = 2403.19995
$ws.Range("L129").Value = 10576.2855
$ws.Range("M129").Value = 2596.80005
$ws.Range("N129").Value = -20576.2855
$ws.Range("H137").Value = 3464.3333
$ws.Range("I137").Value = 3400
$ws.Range("J137").Value = 3496.5
$ws.Range("K137").Value = 10200
$ws.Range("L137").Value = 10489.5
$ws.Range("M137").Value = -5100
$ws.Range("N137").Value = -20689.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3502.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3502.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3502.5
$ws.Range("N80").Value = -5498.5
$ws.Range("H83").Value = 3502.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3502.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 17512.5
$ws.Range("N83").Value = -27496.5
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H97").Value = 1736.6666
$ws.Range("I97").Value = 605
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 605
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = -109
$ws.Range("N97").Value = -4992
$ws.Range("H102").Value = 1895.5454
$ws.Range("I102").Value = 1938.9048
$ws.Range("K102").Value = 1938.9048
$ws.Range("M102").Value = -316.9048
$ws.Range("H108").Value = 52510
$ws.Range("J108").Value = 50020
$ws.Range("L108").Value = 50020
$ws.Range("N108").Value = -57700
$ws.Range("H113").Value = 1813.5
$ws.Range("I113").Value = 1791.3334
$ws.Range("K113").Value = 1791.3334
$ws.Range("M113").Value = 378.6666
$ws.Range("H126").Value = 3381
$ws.Range("I126").Value = 2357.6667
$ws.Range("K126").Value = 7073.000100000001
$ws.Range("M126").Value = -4603.000100000001
$ws.Range("H132").Value = 1955.7
$ws.Range("I132").Value = 1732.125
$ws.Range("J132").Value = 2850
$ws.Range("K132").Value = 5196.375
$ws.Range("L132").Value = 8550
$ws.Range("M132").Value = -2666.375
$ws.Range("N132").Value = -13610
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 77989
$ws.Range("J6").Value = 77989
$ws.Range("L6").Value = 77989
$ws.Range("N6").Value = -78213
$ws.Range("H13").Value = 1120.4
$ws.Range("I13").Value = 129
$ws.Range("J13").Value = 15000
$ws.Range("K13").Value = 129
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 11
$ws.Range("N13").Value = -15280
$ws.Range("H16").Value = 2489.8333
$ws.Range("I16").Value = 1097.7778
$ws.Range("K16").Value = 1097.7778
$ws.Range("M16").Value = -927.7778000000001
$ws.Range("H22").Value = 1107.6428
$ws.Range("I22").Value = 945.1429
$ws.Range("J22").Value = 1270.1428
$ws.Range("K22").Value = 945.1429
$ws.Range("L22").Value = 1270.1428
$ws.Range("M22").Value = -650.1429
$ws.Range("N22").Value = -1860.1428
$ws.Range("H25").Value = 8003.5
$ws.Range("I25").Value = 7
$ws.Range("J25").Value = 16000
$ws.Range("K25").Value = 7
$ws.Range("L25").Value = 16000
$ws.Range("M25").Value = 223
$ws.Range("N25").Value = -16460
$ws.Range("H27").Value = 1107.6428
$ws.Range("I27").Value = 945.1429
$ws.Range("J27").Value = 1270.1428
$ws.Range("K27").Value = 945.1429
$ws.Range("L27").Value = 1270.1428
$ws.Range("M27").Value = -838.1429
$ws.Range("N27").Value = -1484.1428
$ws.Range("H40").Value = 9808110
$ws.Range("I40").Value = 12349778
$ws.Range("K40").Value = 12349778
$ws.Range("M40").Value = -12349642
$ws.Range("H46").Value = 2354.625
$ws.Range("I46").Value = 752
$ws.Range("J46").Value = 3957.25
$ws.Range("K46").Value = 752
$ws.Range("L46").Value = 3957.25
$ws.Range("M46").Value = -564
$ws.Range("N46").Value = -4333.25
$ws.Range("H55").Value = 957.6667
$ws.Range("I55").Value = 775.1111
$ws.Range("K55").Value = 775.1111
$ws.Range("M55").Value = -602.1111
$ws.Range("H61").Value = 4966.1113
$ws.Range("J61").Value = 6500
$ws.Range("L61").Value = 6500
$ws.Range("N61").Value = -6904
$ws.Range("H68").Value = 1953.5
$ws.Range("I68").Value = 1881.7142
$ws.Range("J68").Value = 2054
$ws.Range("K68").Value = 1881.7142
$ws.Range("L68").Value = 2054
$ws.Range("M68").Value = -1132.7142
$ws.Range("N68").Value = -3552
$ws.Range("H71").Value = 1953.5
$ws.Range("I71").Value = 1881.7142
$ws.Range("J71").Value = 2054
$ws.Range("K71").Value = 9408.571
$ws.Range("L71").Value = 10270
$ws.Range("M71").Value = -5664.571
$ws.Range("N71").Value = -17758
$ws.Range("H93").Value = 1128.8823
$ws.Range("I93").Value = 846.7778
$ws.Range("J93").Value = 1446.25
$ws.Range("K93").Value = 846.7778
$ws.Range("L93").Value = 1446.25
$ws.Range("M93").Value = 401.2222
$ws.Range("N93").Value = -3942.25
$ws.Range("H98").Value = 49930
$ws.Range("J98").Value = 49930
$ws.Range("L98").Value = 49930
$ws.Range("N98").Value = -55920
$ws.Range("H100").Value = 1700
$ws.Range("I100").Value = 1700
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1700
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1159
$ws.Range("H113").Value = 4966.1113
$ws.Range("J113").Value = 6500
$ws.Range("L113").Value = 6500
$ws.Range("N113").Value = -10840
$ws.Range("H132").Value = 404282.16
$ws.Range("I132").Value = 717082.5
$ws.Range("J132").Value = 6172.636
$ws.Range("K132").Value = 2151247.5
$ws.Range("L132").Value = 18517.908
$ws.Range("M132").Value = -2148717.5
$ws.Range("N132").Value = -23577.908
$ws.Range("H136").Value = 5306.049
$ws.Range("I136").Value = 4956.2646
$ws.Range("K136").Value = 14868.7938
$ws.Range("M136").Value = -12318.7938
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 60485
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("H100").Value = 1918.7273
$ws.Range("I100").Value = 1903.0625
$ws.Range("J100").Value = 1960.5
$ws.Range("K100").Value = 3806.125
$ws.Range("L100").Value = 3921
$ws.Range("M100").Value = -3265.125
$ws.Range("N100").Value = -5003
$ws.Range("H113").Value = 240.71428
$ws.Range("J113").Value = 309
$ws.Range("L113").Value = 927
$ws.Range("N113").Value = -5267
$ws.Range("H132").Value = 31120.514
$ws.Range("I132").Value = 42915.44
$ws.Range("K132").Value = 128746.32
$ws.Range("M132").Value = -126216.32
$ws.Range("H136").Value = 35413.645
$ws.Range("I136").Value = 2080.125
$ws.Range("J136").Value = 70969.4
$ws.Range("K136").Value = 6240.375
$ws.Range("L136").Value = 212908.2
$ws.Range("M136").Value = -3690.375
$ws.Range("N136").Value = -218008.2
$ws.Range("M99").ClearContents()
